# Stock quantity/value adjustment for CryCompanywiseStockReport
# Updates on-hand quantity (F) and stock value (G = D * F) for affected
# line items, plus the "Sub Total" (B) rows for each vendor group and
# the workbook-wide "Sub Total" / "Grand Total" rows at the bottom.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35: ZOFF Big Cardamom Whole 25Gm
$ws.Range("F35").Value = 187
$ws.Range("G35").Value = 4800.29

# Row 50: ZOFF Green Cardamom Whole 25 GMS
$ws.Range("F50").Value = 14
$ws.Range("G50").Value = 1309.56

# Row 52: ZOFF Indian Bay Leaf 50Gm
$ws.Range("F52").Value = 50
$ws.Range("G52").Value = 820.5

# Row 56: ZOFF Meat Masala 100 GMS
$ws.Range("F56").Value = 17
$ws.Range("G56").Value = 600.4400000000001

# Row 64: subtotal / total
$ws.Range("B64").Value = 168921.37

# Row 95: Bata-Wsp 40Gm Wax Shoe Polish Black Color
$ws.Range("F95").Value = 46
$ws.Range("G95").Value = 1830.8

# Row 96: subtotal / total
$ws.Range("B96").Value = 19135.58

# Row 117: GLT-7O CLOCK PII TWIN BLADE 5S PACK
$ws.Range("F117").Value = 175
$ws.Range("G117").Value = 24561.25

# Row 119: GLT-Gillette Gurad Razor
$ws.Range("F119").Value = 116
$ws.Range("G119").Value = 1918.64

# Row 126: PHP-Ariel Detergnt Powder Matic 1kg FL
$ws.Range("F126").Value = 26
$ws.Range("G126").Value = 5244.46

# Row 136: subtotal / total
$ws.Range("B136").Value = 187520.93

# Row 148: BRILL-Javitri (Mace) 50G
$ws.Range("F148").Value = 30
$ws.Range("G148").Value = 4242.9

# Row 151: BRILL-Laung (Cloves) 100G
$ws.Range("F151").Value = 23
$ws.Range("G151").Value = 2395.45

# Row 155: subtotal / total
$ws.Range("B155").Value = 7144.94

# Row 165: CHO-Cuticura Lavender Mist Talcum Powder 100g
$ws.Range("F165").Value = 11
$ws.Range("G165").Value = 544.28

# Row 177: subtotal / total
$ws.Range("B177").Value = 9980.629999999999

# Row 193: CRO-SStorm2 Wall Fan 450 mm
$ws.Range("F193").Value = 1
$ws.Range("G193").Value = 2518.75

# Row 195: subtotal / total
$ws.Range("B195").Value = 19865.02

# Row 220: DESAI-Ginger Garlic Paste 200g
$ws.Range("F220").Value = 116
$ws.Range("G220").Value = 3599.48

# Row 225: DESAI-Mixed Fruit Jam 500g
$ws.Range("F225").Value = 10
$ws.Range("G225").Value = 812.8

# Row 228: subtotal / total
$ws.Range("B228").Value = 10071.44

# Row 241: GHP-Dish Wash 500 ml
$ws.Range("F241").Value = 118
$ws.Range("G241").Value = 7727.82

# Row 243: GHP-Floor Cleaner ( W ) 1 Ltr
$ws.Range("F243").Value = 79
$ws.Range("G243").Value = 3183.7

# Row 246: GHP-Glamic Dish wash 750 ML
$ws.Range("F246").Value = 88
$ws.Range("G246").Value = 7645.44

# Row 249: GHP-Glamic Disinfactant toilet cleaner 1 Ltr
$ws.Range("F249").Value = 59
$ws.Range("G249").Value = 4433.26

# Row 254: subtotal / total
$ws.Range("B254").Value = 76328.12

# Row 256: GOD-Cin org 100g PO4
$ws.Range("F256").Value = 346
$ws.Range("G256").Value = 39547.8

# Row 257: GOD-Cinthol Original 75 g (Only for South & North East States)
$ws.Range("F257").Value = 2994
$ws.Range("G257").Value = 55389

# Row 259: GOD-GK Flash Combi Pack M89 P48
$ws.Range("F259").Value = 63
$ws.Range("G259").Value = 4341.33

# Row 263: subtotal / total
$ws.Range("B263").Value = 104981.52

# Row 293: HIM-AYURVEDA SANDAL GLOW SOAP 125G IND
$ws.Range("F293").Value = 137
$ws.Range("G293").Value = 4281.25

# Row 323: HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S
$ws.Range("F323").Value = 30
$ws.Range("G323").Value = 2570.4

# Row 325: subtotal / total
$ws.Range("B325").Value = 102034.94

# Row 329: HUL-Bru Inst Poly 50g
$ws.Range("B329").Value = 57077
$ws.Range("D329").Value = 93.08
$ws.Range("E329").Value = 111.2
$ws.Range("F329").Value = 1
$ws.Range("G329").Value = 93.08

# Row 330: HUL-Bru Inst Poly 50g
$ws.Range("B330").Value = 61610
$ws.Range("D330").Value = 102.71
$ws.Range("E330").Value = 122.71
$ws.Range("F330").Value = 62
$ws.Range("G330").Value = 6368.02

# Row 331: HUL-Bru pure coffee GJ 100g
$ws.Range("F331").Value = 4
$ws.Range("G331").Value = 1066.92

# Row 378: HUL-Sunsilk Blk Shn Shp 650Ml
$ws.Range("F378").Value = 35
$ws.Range("G378").Value = 17303.3

# Row 387: HUL-Surf excel qw det liq 500ml
$ws.Range("F387").Value = 5
$ws.Range("G387").Value = 351.55

# Row 391: HUL-taj mahal tea bags 25s
$ws.Range("F391").Value = 8
$ws.Range("G391").Value = 493.12

# Row 395: subtotal / total
$ws.Range("B395").Value = 208389.91

# Row 413: JLM-MBD Baby Dreams WIPES (WITHOUT LID 80S) BUY 1 GET 1 FREE
$ws.Range("F413").Value = 35
$ws.Range("G413").Value = 4355.05

# Row 414: JLM-MBD Shiny Caterpillar Kids Tooth Brush
$ws.Range("F414").Value = 243
$ws.Range("G414").Value = 7834.32

# Row 415: JLM-MBD Shiny Toothbrush Safari
$ws.Range("F415").Value = 260
$ws.Range("G415").Value = 5829.2

# Row 420: subtotal / total
$ws.Range("B420").Value = 41131.29

# Row 433: Garima (Fruit Basket /Mixed Perfume )Premium Incense Stick Zipper Pouch 120 Gm)
$ws.Range("F433").Value = 112
$ws.Range("G433").Value = 4138.4

# Row 436: Nutrigreen Gold 1 Ltr Pouch
$ws.Range("F436").Value = 569
$ws.Range("G436").Value = 79995.71000000001

# Row 438: Nutrigreen Rice Bran Oil 1 Ltr Pouch
$ws.Range("F438").Value = 32
$ws.Range("G438").Value = 4539.52

# Row 439: subtotal / total
$ws.Range("B439").Value = 110591.71

# Row 441: KAR-MYSORE CARBOLIC SOAP - 150 GM
$ws.Range("F441").Value = 31
$ws.Range("G441").Value = 635.1900000000001

# Row 447: subtotal / total
$ws.Range("B447").Value = 982.7

# Row 461: KUS-Cloth Clip Jumbo (For Cloth Stand)- 10 Pcs Pack
$ws.Range("F461").Value = 1
$ws.Range("G461").Value = 68.65000000000001

# Row 474: subtotal / total
$ws.Range("B474").Value = 35225.89

# Row 480: CHUK-Chukde Hing- 50g
$ws.Range("F480").Value = 105
$ws.Range("G480").Value = 10316.25

# Row 493: CHUK-Mirch pwdr 200gm
$ws.Range("F493").Value = 35
$ws.Range("G493").Value = 2115.75

# Row 500: subtotal / total
$ws.Range("B500").Value = 130105.63

# Row 520: CRE-Cremica Chocolate Cream 150Gm
$ws.Range("F520").Value = 159
$ws.Range("G520").Value = 3137.07

# Row 521: CRE-Cremica Classic Crackers(100 +20Gm)120Gm
$ws.Range("F521").Value = 281
$ws.Range("G521").Value = 4616.83

# Row 526: CRE-Cremica Marie Classic 250Gm (6Kg)
$ws.Range("F526").Value = 55
$ws.Range("G526").Value = 1070.3

# Row 534: subtotal / total
$ws.Range("B534").Value = 87750.59

# Row 540: RANGA-Stop - O Bathroom Freshner Power Sprays Engligh Lavender 12ml
$ws.Range("F540").Value = 16
$ws.Range("G540").Value = 1114.4

# Row 541: RANGA-Stop - O Bathroom Freshner Power Sprays Engligh Lemongrass 12ml
$ws.Range("F541").Value = 22
$ws.Range("G541").Value = 1532.3

# Row 544: subtotal / total
$ws.Range("B544").Value = 13678.08

# Row 577: OCT-Octavius Basmati Rice Classic 1 Kg.
$ws.Range("F577").Value = 97
$ws.Range("G577").Value = 14278.4

# Row 581: subtotal / total
$ws.Range("B581").Value = 56866.64

# Row 586: UNB-McvitieS Cream Bourbon (100 Gms)
$ws.Range("F586").Value = 554
$ws.Range("G586").Value = 10996.9

# Row 591: subtotal / total
$ws.Range("B591").Value = 32427.91

# Row 686: KAR-2in1 Biscuit (Fruit+Osmania) - 400gms
$ws.Range("F686").Value = 74
$ws.Range("G686").Value = 9660.700000000001

# Row 689: KAR-Chai Biscuits - 400gms
$ws.Range("F689").Value = 66
$ws.Range("G689").Value = 6266.7

# Row 693: subtotal / total
$ws.Range("B693").Value = 45414.73

# Row 712: SARATHI-Champa 80 Gms Jumbo Economy Box
$ws.Range("F712").Value = 0
$ws.Range("G712").Value = 0

# Row 713: SARATHI-Chandan 20 Sticks Hexagon
$ws.Range("F713").Value = 56
$ws.Range("G713").Value = 888.16

# Row 718: SARATHI-Sandalwood 80 Gms Jumbo Economy Box
$ws.Range("F718").Value = 4
$ws.Range("G718").Value = 132.44

# Row 721: subtotal / total
$ws.Range("B721").Value = 16070.6

# Row 729: SCM-VENSON - COX BRIEF - 95 cm   100 cm
$ws.Range("F729").Value = 5
$ws.Range("G729").Value = 334.65

# Row 737: subtotal / total
$ws.Range("B737").Value = 16917.55

# Row 800: TCP-Chana Dal 1 kg
$ws.Range("F800").Value = 137
$ws.Range("G800").Value = 11173.72

# Row 801: TCP-Coriander powder 200gm
$ws.Range("F801").Value = 116
$ws.Range("G801").Value = 5551.76

# Row 803: TCP-Moong Dal 1kg
$ws.Range("F803").Value = 182
$ws.Range("G803").Value = 23751

# Row 810: TCP-TATA SALT ROCK SALT 1 KG
$ws.Range("F810").Value = 53
$ws.Range("G810").Value = 4806.04

# Row 815: TCP-Tetley Green Tea Teabag Lemon Honey 25s
$ws.Range("F815").Value = 46
$ws.Range("G815").Value = 5098.18

# Row 816: TCP-Toor Dal 1kg
$ws.Range("F816").Value = 300
$ws.Range("G816").Value = 40503

# Row 817: TCP-Urad Dal 1 kg
$ws.Range("F817").Value = 393
$ws.Range("G817").Value = 47439.03

# Row 818: TCP-URAD DAL KALI 10 X 1Kg
$ws.Range("F818").Value = 29
$ws.Range("G818").Value = 3500.59

# Row 819: subtotal / total
$ws.Range("B819").Value = 171950.98

# Row 824: Orgfeed Brown Chana 1 Kg
$ws.Range("F824").Value = 49
$ws.Range("G824").Value = 5331.69

# Row 838: Shankys Tip Top American Popcorn 200 Gm
$ws.Range("F838").Value = 59
$ws.Range("G838").Value = 1950.54

# Row 846: Shankys Tip Top MP Wheat Roasted Dalia 500 Gm
$ws.Range("F846").Value = 118
$ws.Range("G846").Value = 3901.08

# Row 848: Shankys Tip Top Vermicilli 500 Gm
$ws.Range("F848").Value = 119
$ws.Range("G848").Value = 5115.81

# Row 849: Shankys Tip Top Vermicilli Roasted 500 Gm
$ws.Range("F849").Value = 160
$ws.Range("G849").Value = 7982.4

# Row 850: Tip Top Besan 800 g
$ws.Range("F850").Value = 85
$ws.Range("G850").Value = 6851

# Row 854: subtotal / total
$ws.Range("B854").Value = 137798.52

# Row 884: VIP- Conquer 4W Exp Strolly 68 Ast
$ws.Range("F884").Value = 0
$ws.Range("G884").Value = 0

# Row 904: subtotal / total
$ws.Range("B904").Value = 79438.14999999999

# Row 906: VRI-8 AM MASALA OATS 200 gms
$ws.Range("F906").Value = 33
$ws.Range("G906").Value = 1276.77

# Row 911: subtotal / total
$ws.Range("B911").Value = 1336.61

# Row 913: VVD Ayush Cane Jaggery Powder 500Gm
$ws.Range("F913").Value = 268
$ws.Range("G913").Value = 8101.64

# Row 914: VVD Priyam Cold Pressed Groundnut Oil Pouch 1 Ltr
$ws.Range("F914").Value = 2738
$ws.Range("G914").Value = 446595.18

# Row 916: VVD Pure Drop Cold Pressed Gingelly Oil Pouch 500Ml
$ws.Range("F916").Value = 412
$ws.Range("G916").Value = 59595.8

# Row 920: VVD Veda Pancha Deepam Oil Pouch 450Ml
$ws.Range("F920").Value = 153
$ws.Range("G920").Value = 10327.5

# Row 922: subtotal / total
$ws.Range("B922").Value = 662924.4

# Row 928: subtotal / total
$ws.Range("B928").Value = 3961723.51

# Row 929: Note:Rates are Inclusive of Tax
$ws.Range("B929").Value = 3961723.51
